# Apply Tasks-sheet status updates: mark first three tasks as "done" and
# move "in Bearbeitung" to the last task row, then update the selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

$ws.Range("C4").Value = "done"
$ws.Range("C5").Value = "done"
$ws.Range("C6").Value = "done"
$ws.Range("C7").Value = "in Bearbeitung"

$ws.Activate()
$ws.Range("D6").Select()
